# Insert a new weekly price record as row 154 on the single worksheet.
# This pushes the previous rows 154:185 down to 155:186 (data unchanged),
# and grows the sheet's used range from A1:R185 to A1:R186.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 154, shifting rows 154-185 down.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new record.
$ws.Cells.Item(154, 1).Value = 1
$ws.Cells.Item(154, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(154, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(154, 4).Value = 45015
$ws.Cells.Item(154, 5).Value = 15
$ws.Cells.Item(154, 6).Value = 100114001
$ws.Cells.Item(154, 7).Value = 'Papa'
$ws.Cells.Item(154, 8).Value = 'Asterix'
$ws.Cells.Item(154, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(154, 10).Value = 1000
$ws.Cells.Item(154, 11).Value = 15000
$ws.Cells.Item(154, 12).Value = 16000
$ws.Cells.Item(154, 13).Value = 15500
$ws.Cells.Item(154, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(154, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(154, 16).Value = 620
$ws.Cells.Item(154, 17).Value = 25
$ws.Cells.Item(154, 18).Value = 'Hortaliza'
